# Weekly update: a new "Rabanito" (Vega Central Mapocho de Santiago) record for
# the latest date is inserted at row 117, pushing the existing historical rows
# (117-153) down by one (to 118-154). The sheet's used range grows from
# A1:R153 to A1:R154 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 117..153 down to 118..154, inheriting row 116's formatting
# for the freshly inserted row (matches native Excel "Insert Row" behaviour).
$ws.Rows("117:117").Insert()

# Populate the newly inserted row 117 with the new weekly data point.
$ws.Range("A117").Value = 9
$ws.Range("B117").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C117").Value = "Metropolitana"
$ws.Range("D117").Value = 44463
$ws.Range("E117").Value = 13
$ws.Range("F117").Value = 300000001
$ws.Range("G117").Value = "Rabanito"
$ws.Range("H117").Value = "Sin especificar"
$ws.Range("I117").Value = "Primera"
$ws.Range("J117").Value = 7900
$ws.Range("K117").Value = 3500
$ws.Range("L117").Value = 4000
$ws.Range("M117").Value = 3747
$ws.Range("N117").Value = "`$/cien unidades (volumen en unidades)"
$ws.Range("O117").Value = "Provincia de Chacabuco"
$ws.Range("P117").Value = 37
$ws.Range("Q117").Value = 100
$ws.Range("R117").Value = "Hortaliza"
